$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ITC run script "jdccaiicbs.inj" previously used for every non-cleaning
# titration has been replaced by two new, more specific scripts:
#   - ChoderaWaterWater.inj : water-into-water / buffer-into-buffer controls
#   - ChoderaHostGuest.inj  : host/guest (Onesite) titrations
# Rows 2 and 18 (the cleaning-water titrations) keep "water5inj.inj".

$waterWaterRows = 3, 4, 19, 20
foreach ($r in $waterWaterRows) {
    $ws.Cells.Item($r, 4).Value = "ChoderaWaterWater.inj"
}

$hostGuestRows = 5..17
foreach ($r in $hostGuestRows) {
    $ws.Cells.Item($r, 4).Value = "ChoderaHostGuest.inj"
}
